$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-51: Coin, Link, Price, Volume(1h)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.388.32", "  +1.09%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.804.05", "  -0.98%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  -0.04%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "327.63", "  -2.83%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  +0.17%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4450", "  +5.59%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3774", "  +7.04%  "),
    @(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "44.49", "  -2.64%  "),
    @(10, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.149", "  -1.06%  "),
    @(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07498", "  -0.44%  "),
    @(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "22.56", "  -2.37%  "),
    @(13, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  +0.11%  "),
    @(14, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.632", "  +4.19%  "),
    @(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.293", "  -0.49%  "),
    @(16, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.799.37", "  -0.92%  "),
    @(17, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001091", "  -0.20%  "),
    @(18, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06816", "  +1.74%  "),
    @(19, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "80.65", "  -2.83%  "),
    @(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9995", "  -0.08%  "),
    @(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "17.48", "  -0.07%  "),
    @(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.306", "  -1.77%  "),
    @(23, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.361.88", "  +0.78%  "),
    @(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.80", "  -1.39%  "),
    @(25, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.405", "  +0.25%  "),
    @(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.49", "  -2.65%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "154.03", "  -1.61%  "),
    @(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.355", "  -6.84%  "),
    @(29, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.009.13", "  -0.52%  "),
    @(30, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "132.21", "  -1.55%  "),
    @(31, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.255", "  -5.11%  "),
    @(32, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "4.010", "  -1.67%  "),
    @(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.807", "  -4.14%  "),
    @(34, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09319", "  +1.75%  "),
    @(35, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2290", "  +5.19%  "),
    @(36, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "12.11", "  -2.83%  "),
    @(37, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06355", "  -0.04%  "),
    @(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02324", "  -1.48%  "),
    @(39, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6584", "  -1.90%  "),
    @(40, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.158", "  -2.19%  "),
    @(41, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.207", "  -1.18%  "),
    @(42, "WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.458", "  -3.64%  "),
    @(43, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.142", "  -0.50%  "),
    @(44, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9989", "  -0.08%  "),
    @(45, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "13.90", "  -3.45%  "),
    @(46, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.6068", "  -2.01%  "),
    @(47, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.804", "  -1.87%  "),
    @(48, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "128.31", "  -0.39%  "),
    @(49, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "2.030", "  -2.03%  "),
    @(50, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.07092", "  -0.78%  "),
    @(51, "EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.151", "  -3.49%  "),
)

# Ensure Price column (D) is treated as text so values such as "28.388.32"
# or "1.002" are preserved exactly instead of being parsed as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
